# Applies scheduled-runner profit recalculations to the Ixion_Profits sheets.
# For each affected leve row, columns H-N (price/profit figures) are updated
# to reflect refreshed market-board averages.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 49: Going Nowhere Fast | Paralyzing Potion
$ws.Range("H49").Value = 5011.3335
$ws.Range("J49").Value = 7000
$ws.Range("L49").Value = 21000
$ws.Range("N49").Value = -21272

# Row 88: The Grave of Hemlock Groves | Growth Formula Zeta
$ws.Range("H88").Value = 6796.5264
$ws.Range("I88").Value = 1662
$ws.Range("J88").Value = 9791.666999999999
$ws.Range("K88").Value = 1662
$ws.Range("L88").Value = 9791.666999999999
$ws.Range("M88").Value = -1256
$ws.Range("N88").Value = -10603.667

# Row 91: Dappling the Highlands (L) | Growth Formula Zeta
$ws.Range("H91").Value = 6796.5264
$ws.Range("I91").Value = 1662
$ws.Range("J91").Value = 9791.666999999999
$ws.Range("K91").Value = 1662
$ws.Range("L91").Value = 9791.666999999999
$ws.Range("M91").Value = -258
$ws.Range("N91").Value = -12599.667

# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 11905660
$ws.Range("I112").Value = 606.6667
$ws.Range("J112").Value = 15152493
$ws.Range("K112").Value = 1820.0001
$ws.Range("L112").Value = 45457479
$ws.Range("M112").Value = -712.0001
$ws.Range("N112").Value = -45459695

# Row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 5565.769
$ws.Range("I116").Value = 9269.538
$ws.Range("J116").Value = 1862
$ws.Range("K116").Value = 9269.538
$ws.Range("L116").Value = 1862
$ws.Range("M116").Value = -5827.538
$ws.Range("N116").Value = -8746

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 2072.138
$ws.Range("I137").Value = 1503.5454
$ws.Range("J137").Value = 2419.611
$ws.Range("K137").Value = 4510.6362
$ws.Range("L137").Value = 7258.833
$ws.Range("M137").Value = -1960.6362
$ws.Range("N137").Value = -12358.833

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 7108.5
$ws.Range("I32").Value = 7315.85
$ws.Range("J32").Value = 6590.125
$ws.Range("K32").Value = 7315.85
$ws.Range("L32").Value = 6590.125
$ws.Range("M32").Value = -7028.85
$ws.Range("N32").Value = -7164.125

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 5285.4136
$ws.Range("I74").Value = 6677.737
$ws.Range("J74").Value = 2640
$ws.Range("K74").Value = 6677.737
$ws.Range("L74").Value = 2640
$ws.Range("M74").Value = -5803.737
$ws.Range("N74").Value = -4388

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 5285.4136
$ws.Range("I77").Value = 6677.737
$ws.Range("J77").Value = 2640
$ws.Range("K77").Value = 33388.685
$ws.Range("L77").Value = 13200
$ws.Range("M77").Value = -29020.685
$ws.Range("N77").Value = -21936

# Row 88: The Mast Chance | Adamantite Rivets
$ws.Range("H88").Value = 142859220
$ws.Range("I88").Value = 2566
$ws.Range("J88").Value = 250001700
$ws.Range("K88").Value = 2566
$ws.Range("L88").Value = 250001700
$ws.Range("M88").Value = -2160
$ws.Range("N88").Value = -250002512

# Row 91: The Rose and the Riveter (L) | Adamantite Rivets
$ws.Range("H91").Value = 142859220
$ws.Range("I91").Value = 2566
$ws.Range("J91").Value = 250001700
$ws.Range("K91").Value = 2566
$ws.Range("L91").Value = 250001700
$ws.Range("M91").Value = -1162
$ws.Range("N91").Value = -250004508

# Row 139: Backing up My Words | Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 49552.5
$ws.Range("J139").Value = 49552.5
$ws.Range("L139").Value = 49552.5
$ws.Range("N139").Value = -59832.5

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 2355.3076
$ws.Range("I86").Value = 2129
$ws.Range("J86").Value = 3600
$ws.Range("K86").Value = 2129
$ws.Range("L86").Value = 3600
$ws.Range("N86").Value = -5846
$ws.Range("M86").Value = -1006

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 2355.3076
$ws.Range("I89").Value = 2129
$ws.Range("J89").Value = 3600
$ws.Range("K89").Value = 10645
$ws.Range("L89").Value = 18000
$ws.Range("N89").Value = -29232
$ws.Range("M89").Value = -5029

$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face | Fermented Butter
$ws.Range("H68").Value = 2934.4204
$ws.Range("I68").Value = 3917.7273
$ws.Range("J68").Value = 2033.0555
$ws.Range("K68").Value = 11753.1819
$ws.Range("L68").Value = 6099.166499999999
$ws.Range("M68").Value = -10942.1819
$ws.Range("N68").Value = -7721.166499999999

# Row 71: No Margarine of Error (L) | Fermented Butter
$ws.Range("H71").Value = 2934.4204
$ws.Range("I71").Value = 3917.7273
$ws.Range("J71").Value = 2033.0555
$ws.Range("K71").Value = 35259.5457
$ws.Range("L71").Value = 18297.4995
$ws.Range("M71").Value = -31203.5457
$ws.Range("N71").Value = -26409.4995

$ws = $wb.Worksheets.Item("GSM")
# Row 112: Gentleman Donor | Diaspore Bracelet of Slaying
$ws.Range("H112").Value = 37200
$ws.Range("J112").Value = 37200
$ws.Range("L112").Value = 37200
$ws.Range("N112").Value = -39416

# Row 114: Hot Rod | Bluespirit Rod
$ws.Range("H114").Value = 47500
$ws.Range("J114").Value = 47500
$ws.Range("L114").Value = 47500
$ws.Range("N114").Value = -56178

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 6482995
$ws.Range("I122").Value = 7203183
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 21609549
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -21607099
$ws.Range("N122").Value = -8800

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 8548542
$ws.Range("I22").Value = 22223140
$ws.Range("J22").Value = 1918.625
$ws.Range("K22").Value = 22223140
$ws.Range("L22").Value = 1918.625
$ws.Range("M22").Value = -22222845
$ws.Range("N22").Value = -2508.625

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 8548542
$ws.Range("I27").Value = 22223140
$ws.Range("J27").Value = 1918.625
$ws.Range("K27").Value = 22223140
$ws.Range("L27").Value = 1918.625
$ws.Range("M27").Value = -22223033
$ws.Range("N27").Value = -2132.625

# Row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws.Range("H55").Value = 15151823
$ws.Range("I55").Value = 238.35294
$ws.Range("J55").Value = 31250382
$ws.Range("K55").Value = 238.35294
$ws.Range("L55").Value = 31250382
$ws.Range("M55").Value = -65.35293999999999
$ws.Range("N55").Value = -31250728

# Row 114: A Heady Endeavor | Atrociraptorskin Headgear of Scouting
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 3139327.5
$ws.Range("I122").Value = 4208472.5
$ws.Range("J122").Value = 1119831.6
$ws.Range("K122").Value = 12625417.5
$ws.Range("L122").Value = 3359494.8
$ws.Range("M122").Value = -12622967.5
$ws.Range("N122").Value = -3364394.8

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display | Ruby Cotton Cloth
$ws.Range("H96").Value = 1964.85
$ws.Range("I96").Value = 1338
$ws.Range("J96").Value = 2477.7273
$ws.Range("K96").Value = 1338
$ws.Range("L96").Value = 2477.7273
$ws.Range("M96").Value = 35
$ws.Range("N96").Value = -5223.7273

# Row 107: Flax Wax | Bright Linen Yarn
$ws.Range("H107").Value = 125000650
$ws.Range("I107").Value = 142857780
$ws.Range("K107").Value = 428573340
$ws.Range("M107").Value = -428571420

# Row 123: Helping Handwear | Fingerless Darkhempen Gloves of Healing
$ws.Range("H123").Value = 26738.6
$ws.Range("J123").Value = 26738.6
$ws.Range("L123").Value = 26738.6
$ws.Range("N123").Value = -36538.6

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 1192.3334
$ws.Range("I126").Value = 829.0833
$ws.Range("K126").Value = 2487.2499
$ws.Range("M126").Value = -17.2498999999998

# Row 138: Halfgloves, Full Effort | Rroneek Serge Halfgloves of Healing
$ws.Range("H138").Value = 23124.5
$ws.Range("J138").Value = 23124.5
$ws.Range("L138").Value = 23124.5
$ws.Range("N138").Value = -33404.5
